# The document's Pearson / BTEC logo images (inline pictures living in the
# footers and header) had their display names re-shuffled:
#   - Footer 1's logo:  image1.png -> image2.png
#   - Footer 2's logo:  image1.png -> image2.png
#   - Header 2's logo:  image2.jpg -> image1.jpg
#
# These are purely the "Name" of the InlineShape (the <wp:docPr name="…"/>
# that backs the picture), so rename each shape through the Word object
# model via Sections -> Headers/Footers -> Range.InlineShapes.

$d = $word.ActiveDocument

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections($s)

    # --- Footers -------------------------------------------------------
    for ($i = 1; $i -le 3; $i++) {
        $footer = $section.Footers($i)
        if ($footer.Exists) {
            $shapes = $footer.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shape = $shapes($j)
                if ($shape.Name -eq "image1.png") {
                    $shape.Name = "image2.png"
                }
            }
        }
    }

    # --- Headers ---------------------------------------------------------
    for ($i = 1; $i -le 3; $i++) {
        $header = $section.Headers($i)
        if ($header.Exists) {
            $shapes = $header.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shape = $shapes($j)
                if ($shape.Name -eq "image2.jpg") {
                    $shape.Name = "image1.jpg"
                }
            }
        }
    }
}
